# Select the worksheet that holds the roster data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Insert a brand new column before column C. This shifts the existing
# C:G columns (one position) to the right, becoming D:H, and creates a
# fresh, empty column C that will be used to hold a "group" indicator.
$ws.Columns("C").Insert()

# Populate the new column C with a group label for every data row.
# Rows 2-5 belong to "group1" and rows 6-15 belong to "group2", except
# row 10 which has no underlying data and is therefore left blank.
$groups = @{
    2  = "group1"
    3  = "group1"
    4  = "group1"
    5  = "group1"
    6  = "group2"
    7  = "group2"
    8  = "group2"
    9  = "group2"
    11 = "group2"
    12 = "group2"
    13 = "group2"
    14 = "group2"
    15 = "group2"
}

foreach ($row in $groups.Keys) {
    $ws.Cells.Item($row, 3).Value = $groups[$row]
}

Write-Host "Inserted group column and populated group values."
